$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header labels
$ws.Range("I1").Value = "Electric Choice ID"
$ws.Range("J1").Value = "Electric Rate Code"
$ws.Range("L1").Value = "Electric Usage (kWh)"

# Add new headers for gas columns
$ws.Range("M1").Value = "Gas Choice ID"
$ws.Range("N1").Value = "Gas Rate Code"
$ws.Range("O1").Value = "Gas Usage (therms)"

# Match the header style (bold, bordered, centered) used by the rest of row 1
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in K2 (previously blank) with N/A
$ws.Range("K2").Value = "N/A"

# Add new data for gas columns in row 2
$ws.Range("M2").Value = "N/A"
$ws.Range("N2").Value = "N/A"
$ws.Range("O2").Value = "N/A"
